$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 395
$ws.Range("B4").Value = 404
$ws.Range("A11").Value = "13.02.2024- Otistics Kel'el Ware karşılığında Out of Po'ya 2 Dolar vermiştir. (395-404)"

$ws.Range("B17").Select()
